# Apply newly added iAuthor TC rows (rows 2-15, columns A-V)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 14,22
$data[0,0] = 1
$data[0,1] = 1
$data[0,2] = "MCQ"
$data[0,3] = "Nill"
$data[0,4] = "Nill"
$data[0,5] = "Nill"
$data[0,6] = "C"
$data[0,7] = 8.5
$data[0,8] = -0.06
$data[0,9] = 27.8
$data[0,10] = -0.03
$data[0,11] = 25.1
$data[0,12] = 0.21
$data[0,13] = 21.3
$data[0,14] = -0.07
$data[0,15] = 17.3
$data[0,16] = -0.09
$data[0,17] = "Nill"
$data[0,18] = "Nill"
$data[0,19] = 0.21
$data[0,20] = 2.25
$data[0,21] = "Low"
$data[1,0] = 2
$data[1,1] = 2
$data[1,2] = "VSAQ"
$data[1,3] = "Nill"
$data[1,4] = "Nill"
$data[1,5] = "Nill"
$data[1,6] = "Nill"
$data[1,7] = "Nill"
$data[1,8] = "Nill"
$data[1,9] = "Nill"
$data[1,10] = "Nill"
$data[1,11] = "Nill"
$data[1,12] = "Nill"
$data[1,13] = "Nill"
$data[1,14] = "Nill"
$data[1,15] = "Nill"
$data[1,16] = "Nill"
$data[1,17] = 55.0
$data[1,18] = "Nill"
$data[1,19] = 0.45
$data[1,20] = -0.9
$data[1,21] = "Low"
$data[2,0] = 3
$data[2,1] = 3
$data[2,2] = "ISAWE"
$data[2,3] = 1
$data[2,4] = 1
$data[2,5] = "Nill"
$data[2,6] = "Nill"
$data[2,7] = "Nill"
$data[2,8] = "Nill"
$data[2,9] = "Nill"
$data[2,10] = "Nill"
$data[2,11] = "Nill"
$data[2,12] = "Nill"
$data[2,13] = "Nill"
$data[2,14] = "Nill"
$data[2,15] = "Nill"
$data[2,16] = "Nill"
$data[2,17] = 35.67
$data[2,18] = "Nill"
$data[2,19] = 0.32
$data[2,20] = -0.9
$data[2,21] = "Low"
$data[3,0] = 4
$data[3,1] = 3
$data[3,2] = "ISAWE"
$data[3,3] = 1
$data[3,4] = 2
$data[3,5] = "Nill"
$data[3,6] = "Nill"
$data[3,7] = "Nill"
$data[3,8] = "Nill"
$data[3,9] = "Nill"
$data[3,10] = "Nill"
$data[3,11] = "Nill"
$data[3,12] = "Nill"
$data[3,13] = "Nill"
$data[3,14] = "Nill"
$data[3,15] = "Nill"
$data[3,16] = "Nill"
$data[3,17] = 65.32
$data[3,18] = "Nill"
$data[3,19] = 0.21
$data[3,20] = -0.9
$data[3,21] = "Low"
$data[4,0] = 5
$data[4,1] = 3
$data[4,2] = "ISAWE"
$data[4,3] = 2
$data[4,4] = 1
$data[4,5] = "Nill"
$data[4,6] = "Nill"
$data[4,7] = "Nill"
$data[4,8] = "Nill"
$data[4,9] = "Nill"
$data[4,10] = "Nill"
$data[4,11] = "Nill"
$data[4,12] = "Nill"
$data[4,13] = "Nill"
$data[4,14] = "Nill"
$data[4,15] = "Nill"
$data[4,16] = "Nill"
$data[4,17] = 54.12
$data[4,18] = "Nill"
$data[4,19] = 0.37
$data[4,20] = -0.9
$data[4,21] = "High"
$data[5,0] = 6
$data[5,1] = 3
$data[5,2] = "ISAWE"
$data[5,3] = 2
$data[5,4] = 2
$data[5,5] = "Nill"
$data[5,6] = "Nill"
$data[5,7] = "Nill"
$data[5,8] = "Nill"
$data[5,9] = "Nill"
$data[5,10] = "Nill"
$data[5,11] = "Nill"
$data[5,12] = "Nill"
$data[5,13] = "Nill"
$data[5,14] = "Nill"
$data[5,15] = "Nill"
$data[5,16] = "Nill"
$data[5,17] = 38.05
$data[5,18] = "Nill"
$data[5,19] = 0.54
$data[5,20] = -0.9
$data[5,21] = "High"
$data[6,0] = 7
$data[6,1] = 4
$data[6,2] = "Type X"
$data[6,3] = "Nill"
$data[6,4] = "Nill"
$data[6,5] = "A"
$data[6,6] = "F"
$data[6,7] = 59.9
$data[6,8] = -0.25
$data[6,9] = "Nill"
$data[6,10] = "Nill"
$data[6,11] = "Nill"
$data[6,12] = "Nill"
$data[6,13] = "Nill"
$data[6,14] = "Nill"
$data[6,15] = "Nill"
$data[6,16] = "Nill"
$data[6,17] = "Nill"
$data[6,18] = "Nill"
$data[6,19] = 0.26
$data[6,20] = 1.54
$data[6,21] = "Medium"
$data[7,0] = 8
$data[7,1] = 4
$data[7,2] = "Type X"
$data[7,3] = "Nill"
$data[7,4] = "Nill"
$data[7,5] = "B"
$data[7,6] = "T"
$data[7,7] = "Nill"
$data[7,8] = "Nill"
$data[7,9] = 12.6
$data[7,10] = -0.28
$data[7,11] = "Nill"
$data[7,12] = "Nill"
$data[7,13] = "Nill"
$data[7,14] = "Nill"
$data[7,15] = "Nill"
$data[7,16] = "Nill"
$data[7,17] = "Nill"
$data[7,18] = "Nill"
$data[7,19] = 0.28
$data[7,20] = -0.9
$data[7,21] = "Medium"
$data[8,0] = 9
$data[8,1] = 4
$data[8,2] = "Type X"
$data[8,3] = "Nill"
$data[8,4] = "Nill"
$data[8,5] = "C"
$data[8,6] = "T"
$data[8,7] = "Nill"
$data[8,8] = "Nill"
$data[8,9] = "Nill"
$data[8,10] = "Nill"
$data[8,11] = 39.8
$data[8,12] = 0.26
$data[8,13] = "Nill"
$data[8,14] = "Nill"
$data[8,15] = "Nill"
$data[8,16] = "Nill"
$data[8,17] = "Nill"
$data[8,18] = "Nill"
$data[8,19] = 0.28
$data[8,20] = -0.9
$data[8,21] = "Low"
$data[9,0] = 10
$data[9,1] = 4
$data[9,2] = "Type X"
$data[9,3] = "Nill"
$data[9,4] = "Nill"
$data[9,5] = "D"
$data[9,6] = "T"
$data[9,7] = "Nill"
$data[9,8] = "Nill"
$data[9,9] = "Nill"
$data[9,10] = "Nill"
$data[9,11] = "Nill"
$data[9,12] = "Nill"
$data[9,13] = 12.6
$data[9,14] = -0.28
$data[9,15] = "Nill"
$data[9,16] = "Nill"
$data[9,17] = "Nill"
$data[9,18] = "Nill"
$data[9,19] = 0.28
$data[9,20] = -0.9
$data[9,21] = "High"
$data[10,0] = 11
$data[10,1] = 4
$data[10,2] = "Type X"
$data[10,3] = "Nill"
$data[10,4] = "Nill"
$data[10,5] = "E"
$data[10,6] = "T"
$data[10,7] = "Nill"
$data[10,8] = "Nill"
$data[10,9] = "Nill"
$data[10,10] = "Nill"
$data[10,11] = "Nill"
$data[10,12] = "Nill"
$data[10,13] = "Nill"
$data[10,14] = "Nill"
$data[10,15] = 87.4
$data[10,16] = 0.28
$data[10,17] = "Nill"
$data[10,18] = "Nill"
$data[10,19] = 0.28
$data[10,20] = -0.9
$data[10,21] = "Low"
$data[11,0] = 12
$data[11,1] = 5
$data[11,2] = "Type B"
$data[11,3] = "Nill"
$data[11,4] = "Nill"
$data[11,5] = "Nill"
$data[11,6] = "C"
$data[11,7] = 8.5
$data[11,8] = -0.06
$data[11,9] = 27.8
$data[11,10] = -0.03
$data[11,11] = 25.1
$data[11,12] = 0.21
$data[11,13] = 21.3
$data[11,14] = -0.07
$data[11,15] = 17.3
$data[11,16] = -0.09
$data[11,17] = "Nill"
$data[11,18] = "Nill"
$data[11,19] = 0.21
$data[11,20] = 2.25
$data[11,21] = "Low"
$data[12,0] = 13
$data[12,1] = 6
$data[12,2] = "SAQ"
$data[12,3] = "Nill"
$data[12,4] = "Nill"
$data[12,5] = "Nill"
$data[12,6] = "Nill"
$data[12,7] = "Nill"
$data[12,8] = "Nill"
$data[12,9] = "Nill"
$data[12,10] = "Nill"
$data[12,11] = "Nill"
$data[12,12] = "Nill"
$data[12,13] = "Nill"
$data[12,14] = "Nill"
$data[12,15] = "Nill"
$data[12,16] = "Nill"
$data[12,17] = 55.0
$data[12,18] = "Nill"
$data[12,19] = 0.45
$data[12,20] = -0.9
$data[12,21] = "Low"
$data[13,0] = 14
$data[13,1] = 7
$data[13,2] = "SJT"
$data[13,3] = "Nill"
$data[13,4] = "Nill"
$data[13,5] = "Nill"
$data[13,6] = "Nill"
$data[13,7] = "Nill"
$data[13,8] = "Nill"
$data[13,9] = "Nill"
$data[13,10] = "Nill"
$data[13,11] = "Nill"
$data[13,12] = "Nill"
$data[13,13] = "Nill"
$data[13,14] = "Nill"
$data[13,15] = "Nill"
$data[13,16] = "Nill"
$data[13,17] = 55.0
$data[13,18] = "Nill"
$data[13,19] = 0.45
$data[13,20] = -0.9
$data[13,21] = "Low"

$ws.Range("A2:V15").Value = $data

